# Atualização de bases das ligas, do dia: 16-05-2024 às 23:38
#
# The underlying data rows were re-sorted; for six pairs of adjacent rows the
# match record (every column except the leading sequential-id column A) was
# swapped between the two rows. Column A keeps its original sequential value
# per row; columns B (match id) through AB (odds) move with the record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (row1, row2) pairs whose B:AB contents get swapped
$pairs = @(
    @(29, 30),
    @(36, 37),
    @(49, 50),
    @(76, 77),
    @(87, 88),
    @(177, 178)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AB$r1")
    $range2 = $ws.Range("B$r2`:AB$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
